$p = $ppt.ActivePresentation
try {
    $p.NotesMaster = $null
    Write-Output "set ok"
} catch {
    Write-Output "ERROR: $_"
}
